$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value2 = 1699.8928
$ws.Range("I137").Value2 = 1247
$ws.Range("J137").Value2 = 2303.75
$ws.Range("K137").Value2 = 3741
$ws.Range("L137").Value2 = 6911.25
$ws.Range("M137").Value2 = -1191
$ws.Range("N137").Value2 = -12011.25
$ws.Range("H138").Value2 = 2445.55
$ws.Range("J138").Value2 = 3117.5386
$ws.Range("L138").Value2 = 9352.6158
$ws.Range("N138").Value2 = -19632.6158

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value2 = 5055.5
$ws.Range("I45").Value2 = 1499
$ws.Range("K45").Value2 = 1499
$ws.Range("M45").Value2 = -1122
$ws.Range("H61").Value2 = 731.5
$ws.Range("I61").Value2 = 637.8
$ws.Range("K61").Value2 = 637.8
$ws.Range("M61").Value2 = -425.8
$ws.Range("H63").Value2 = 9339.714
$ws.Range("I63").Value2 = 15224.4375
$ws.Range("K63").Value2 = 15224.4375
$ws.Range("M63").Value2 = -14538.4375
$ws.Range("H66").Value2 = 9339.714
$ws.Range("I66").Value2 = 15224.4375
$ws.Range("K66").Value2 = 76122.1875
$ws.Range("M66").Value2 = -72690.1875
$ws.Range("H74").Value2 = 1423.5
$ws.Range("I74").Value2 = 1397.7391
$ws.Range("K74").Value2 = 1397.7391
$ws.Range("M74").Value2 = -523.7391
$ws.Range("H77").Value2 = 1423.5
$ws.Range("I77").Value2 = 1397.7391
$ws.Range("K77").Value2 = 6988.6955
$ws.Range("M77").Value2 = -2620.6955
$ws.Range("H97").Value2 = 2271.5454
$ws.Range("I97").Value2 = 1973.2778
$ws.Range("K97").Value2 = 1973.2778
$ws.Range("M97").Value2 = -1477.2778
$ws.Range("H125").Value2 = 50715
$ws.Range("J125").Value2 = 50715
$ws.Range("L125").Value2 = 50715
$ws.Range("N125").Value2 = -60555
$ws.Range("H132").Value2 = 1454.909
$ws.Range("I132").Value2 = 1290
$ws.Range("K132").Value2 = 3870
$ws.Range("M132").Value2 = -1340
$ws.Range("H136").Value2 = 731.5
$ws.Range("I136").Value2 = 637.8
$ws.Range("K136").Value2 = 1913.4
$ws.Range("M136").Value2 = 636.6000000000001

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value2 = 9219.091
$ws.Range("I20").Value2 = 10223.333
$ws.Range("J20").Value2 = 4700
$ws.Range("K20").Value2 = 10223.333
$ws.Range("L20").Value2 = 4700
$ws.Range("M20").Value2 = -9976.333000000001
$ws.Range("N20").Value2 = -5194
$ws.Range("H81").Value2 = 39999.5
$ws.Range("J81").Value2 = 39999.5
$ws.Range("L81").Value2 = 39999.5
$ws.Range("N81").Value2 = -42121.5
$ws.Range("H84").Value2 = 39999.5
$ws.Range("J84").Value2 = 39999.5
$ws.Range("L84").Value2 = 119998.5
$ws.Range("N84").Value2 = -130606.5
$ws.Range("H103").Value2 = 50954.332
$ws.Range("J103").Value2 = 50954.332
$ws.Range("L103").Value2 = 50954.332
$ws.Range("N103").Value2 = -53298.332
$ws.Range("H134").Value2 = 2393.4666
$ws.Range("I134").Value2 = 2559.963
$ws.Range("K134").Value2 = 7679.889000000001
$ws.Range("M134").Value2 = -5144.889000000001

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value2 = 6579128
$ws.Range("I7").Value2 = 7812683
$ws.Range("J7").Value2 = 166.33333
$ws.Range("K7").Value2 = 7812683
$ws.Range("L7").Value2 = 166.33333
$ws.Range("M7").Value2 = -7812570
$ws.Range("N7").Value2 = -392.33333
$ws.Range("H62").Value2 = 40918.75
$ws.Range("I62").Value2 = 3592.3076
$ws.Range("J62").Value2 = 202666.67
$ws.Range("K62").Value2 = 3592.3076
$ws.Range("L62").Value2 = 202666.67
$ws.Range("M62").Value2 = -2968.3076
$ws.Range("N62").Value2 = -203914.67
$ws.Range("H65").Value2 = 40918.75
$ws.Range("I65").Value2 = 3592.3076
$ws.Range("J65").Value2 = 202666.67
$ws.Range("K65").Value2 = 17961.538
$ws.Range("L65").Value2 = 1013333.35
$ws.Range("M65").Value2 = -14841.538
$ws.Range("N65").Value2 = -1019573.35
$ws.Range("H132").Value2 = 1194.2
$ws.Range("I132").Value2 = 1105.7778
$ws.Range("J132").Value2 = 1990
$ws.Range("K132").Value2 = 3317.3334
$ws.Range("L132").Value2 = 5970
$ws.Range("M132").Value2 = -787.3334000000004
$ws.Range("N132").Value2 = -11030
$ws.Range("H134").Value2 = 2603.3794
$ws.Range("I134").Value2 = 2463
$ws.Range("J134").Value2 = 4498.5
$ws.Range("K134").Value2 = 7389
$ws.Range("L134").Value2 = 13495.5
$ws.Range("M134").Value2 = -4854
$ws.Range("N134").Value2 = -18565.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value2 = 216076.14
$ws.Range("I4").Value2 = 258.5
$ws.Range("J4").Value2 = 503833
$ws.Range("K4").Value2 = 775.5
$ws.Range("L4").Value2 = 1511499
$ws.Range("M4").Value2 = -663.5
$ws.Range("N4").Value2 = -1511723
$ws.Range("H5").Value2 = 275.25
$ws.Range("I5").Value2 = 275.25
$ws.Range("J5").Value2 = 0
$ws.Range("K5").Value2 = 825.75
$ws.Range("L5").Value2 = 0
$ws.Range("M5").Value2 = -713.75
$ws.Range("H23").Value2 = 531.7
$ws.Range("J23").Value2 = 893
$ws.Range("L23").Value2 = 2679
$ws.Range("N23").Value2 = -3149
$ws.Range("H33").Value2 = 45
$ws.Range("I33").Value2 = 45
$ws.Range("K33").Value2 = 270
$ws.Range("M33").Value2 = 13
$ws.Range("H61").Value2 = 3224
$ws.Range("I61").Value2 = 3788.8
$ws.Range("J61").Value2 = 400
$ws.Range("K61").Value2 = 11366.4
$ws.Range("L61").Value2 = 1200
$ws.Range("M61").Value2 = -11151.4
$ws.Range("N61").Value2 = -1630
$ws.Range("H97").Value2 = 6501.1177
$ws.Range("I97").Value2 = 596.25
$ws.Range("J97").Value2 = 11749.889
$ws.Range("K97").Value2 = 1788.75
$ws.Range("L97").Value2 = 35249.667
$ws.Range("M97").Value2 = -1292.75
$ws.Range("N97").Value2 = -36241.667
$ws.Range("H107").Value2 = 1690.6666
$ws.Range("I107").Value2 = 850
$ws.Range("J107").Value2 = 1858.8
$ws.Range("K107").Value2 = 2550
$ws.Range("L107").Value2 = 5576.4
$ws.Range("M107").Value2 = -630
$ws.Range("N107").Value2 = -9416.4
$ws.Range("H112").Value2 = 3229.4
$ws.Range("J112").Value2 = 4030
$ws.Range("L112").Value2 = 12090
$ws.Range("N112").Value2 = -14306
$ws.Range("H117").Value2 = 1514.2222
$ws.Range("J117").Value2 = 974.75
$ws.Range("L117").Value2 = 2924.25
$ws.Range("N117").Value2 = -9808.25
$ws.Range("H121").Value2 = 2007704.8
$ws.Range("I121").Value2 = 145168.28
$ws.Range("J121").Value2 = 3010609
$ws.Range("K121").Value2 = 435504.84
$ws.Range("L121").Value2 = 9031827
$ws.Range("M121").Value2 = -434194.84
$ws.Range("N121").Value2 = -9034447
$ws.Range("H135").Value2 = 275.25
$ws.Range("I135").Value2 = 275.25
$ws.Range("J135").Value2 = 0
$ws.Range("K135").Value2 = 2477.25
$ws.Range("L135").Value2 = 0
$ws.Range("M135").Value2 = 57.75
$ws.Range("H140").Value2 = 4878.727
$ws.Range("I140").Value2 = 5583.375
$ws.Range("J140").Value2 = 2999.6667
$ws.Range("K140").Value2 = 16750.125
$ws.Range("L140").Value2 = 8999.000100000001
$ws.Range("M140").Value2 = -11570.125
$ws.Range("N140").Value2 = -19359.0001
$ws.Range("N5").ClearContents()
$ws.Range("N135").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H68").Value2 = 40268
$ws.Range("I68").Value2 = 40268
$ws.Range("K68").Value2 = 40268
$ws.Range("M68").Value2 = -39457
$ws.Range("H71").Value2 = 40268
$ws.Range("I71").Value2 = 40268
$ws.Range("K71").Value2 = 120804
$ws.Range("M71").Value2 = -116748
$ws.Range("H107").Value2 = 1135.8334
$ws.Range("I107").Value2 = 1165
$ws.Range("K107").Value2 = 1165
$ws.Range("M107").Value2 = 755
$ws.Range("H113").Value2 = 2997.25
$ws.Range("I113").Value2 = 2994.5
$ws.Range("K113").Value2 = 2994.5
$ws.Range("M113").Value2 = -824.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value2 = 646.7406999999999
$ws.Range("I16").Value2 = 680.92
$ws.Range("K16").Value2 = 680.92
$ws.Range("M16").Value2 = -510.92
$ws.Range("H24").Value2 = 21996
$ws.Range("I24").Value2 = 21996
$ws.Range("K24").Value2 = 21996
$ws.Range("M24").Value2 = -21653
$ws.Range("H100").Value2 = 1505002
$ws.Range("I100").Value2 = 3000000
$ws.Range("K100").Value2 = 3000000
$ws.Range("M100").Value2 = -2999459
$ws.Range("H132").Value2 = 7107.357
$ws.Range("I132").Value2 = 6624.1113
$ws.Range("J132").Value2 = 7977.2
$ws.Range("K132").Value2 = 19872.3339
$ws.Range("L132").Value2 = 23931.6
$ws.Range("M132").Value2 = -17342.3339
$ws.Range("N132").Value2 = -28991.6
$ws.Range("H136").Value2 = 83337784
$ws.Range("I136").Value2 = 4300.8887
$ws.Range("K136").Value2 = 12902.6661
$ws.Range("M136").Value2 = -10352.6661

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H30").Value2 = 24478.666
$ws.Range("J30").Value2 = 24974.4
$ws.Range("L30").Value2 = 24974.4
$ws.Range("N30").Value2 = -25188.4
$ws.Range("H101").Value2 = 15150.667
$ws.Range("J101").Value2 = 15150.667
$ws.Range("L101").Value2 = 15150.667
$ws.Range("N101").Value2 = -21640.667
$ws.Range("H107").Value2 = 1811.1364
$ws.Range("J107").Value2 = 3292.375
$ws.Range("L107").Value2 = 9877.125
$ws.Range("N107").Value2 = -13717.125
$ws.Range("H132").Value2 = 3971.524
$ws.Range("I132").Value2 = 2435.3333
$ws.Range("J132").Value2 = 5123.6665
$ws.Range("K132").Value2 = 7305.999899999999
$ws.Range("L132").Value2 = 15370.9995
$ws.Range("M132").Value2 = -4775.999899999999
$ws.Range("N132").Value2 = -20430.9995
$ws.Range("H136").Value2 = 1114
$ws.Range("I136").Value2 = 1051
$ws.Range("K136").Value2 = 3153
$ws.Range("M136").Value2 = -603
